# Property.xlsx - "login coded updated, login tbl"
# Adds a "Login" fields list (column O) and a "Verification" fields list
# (column P) to Sheet1, to the right of the existing Payment/Tenant table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (bold + italic, matching the other headers in row 2) ----
$ws.Range("O2").Value = "Login"
$ws.Range("P2").Value = "Verification"
$ws.Range("O2:P2").Font.Bold = $true
$ws.Range("O2:P2").Font.Italic = $true

# ---- Column O: Login fields (rows 2-11 first, then 12-16 later) ----
$ws.Range("O3").Value = "Id"
$ws.Range("O4").Value = "UserName"
$ws.Range("O5").Value = "Password"
$ws.Range("O6").Value = "HouseOwnerId"
$ws.Range("O7").Value = "TenantId"
$ws.Range("O8").Value = "PhoneNumber"
$ws.Range("O9").Value = "EmailID"
$ws.Range("O10").Value = "PhoneNumberVerified"
$ws.Range("O11").Value = "EmailIdVerified"

# ---- Column P: Verification fields ----
$ws.Range("P3").Value = "Id"
$ws.Range("P4").Value = "PhoneNumber"
$ws.Range("P5").Value = "PhoneNumberCode"
$ws.Range("P6").Value = "CodeSendDate"
$ws.Range("P8").Value = "ExipryTime"
$ws.Range("P9").Value = "Status"
$ws.Range("P10").Value = "Email"
$ws.Range("P11").Value = "EMailCode"
$ws.Range("P7").Value = "CodeSendTime"
$ws.Range("P12").Value = "LoginID"

# ---- Column O continued: remaining Login fields ----
$ws.Range("O12").Value = "ReVerificationTime"
$ws.Range("O13").Value = "PhoneNumberVerifiedDate"
$ws.Range("O14").Value = "EmailIdVerifiedDate"
$ws.Range("O15").Value = "MandatoryVerification"
$ws.Range("O16").Value = "ReVerification"

# ---- Column widths for the two new columns ----
$ws.Columns.Item(15).ColumnWidth = 25.17
$ws.Columns.Item(16).ColumnWidth = 18.02

# ---- View state: scroll/selection ends on the last filled cell ----
$ws.Range("O16").Select()
